# Auto-generated Excel COM-interop script to apply Sophia_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 877.2941
$ws.Range("I92").Value = 1007.1818
$ws.Range("J92").Value = 639.1667
$ws.Range("K92").Value = 1007.1818
$ws.Range("L92").Value = 639.1667
$ws.Range("M92").Value = 240.8182
$ws.Range("N92").Value = -3135.1667
$ws.Range("H96").Value = 366.41666
$ws.Range("I96").Value = 222.9
$ws.Range("K96").Value = 668.7
$ws.Range("M96").Value = 704.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1900
$ws.Range("I61").Value = 1900
$ws.Range("K61").Value = 1900
$ws.Range("M61").Value = -1688
$ws.Range("H122").Value = 3345.3076
$ws.Range("I122").Value = 2953.4546
$ws.Range("J122").Value = 5500.5
$ws.Range("K122").Value = 8860.363799999999
$ws.Range("L122").Value = 16501.5
$ws.Range("M122").Value = -6410.363799999999
$ws.Range("N122").Value = -21401.5
$ws.Range("H136").Value = 1900
$ws.Range("I136").Value = 1900
$ws.Range("K136").Value = 5700
$ws.Range("M136").Value = -3150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3597.4167
$ws.Range("I94").Value = 2833.875
$ws.Range("J94").Value = 5124.5
$ws.Range("K94").Value = 2833.875
$ws.Range("L94").Value = 5124.5
$ws.Range("M94").Value = -2382.875
$ws.Range("N94").Value = -6026.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 162.90909
$ws.Range("I7").Value = 132.22223
$ws.Range("J7").Value = 301
$ws.Range("K7").Value = 132.22223
$ws.Range("L7").Value = 301
$ws.Range("M7").Value = -19.22223
$ws.Range("N7").Value = -527
$ws.Range("H134").Value = 4943.353
$ws.Range("I134").Value = 4673.8184
$ws.Range("K134").Value = 14021.4552
$ws.Range("M134").Value = -11486.4552
$ws.Range("H141").Value = 50998.5
$ws.Range("J141").Value = 50998.5
$ws.Range("L141").Value = 50998.5
$ws.Range("N141").Value = -61358.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 389.375
$ws.Range("I11").Value = 162.5
$ws.Range("J11").Value = 1070
$ws.Range("K11").Value = 487.5
$ws.Range("L11").Value = 3210
$ws.Range("M11").Value = -347.5
$ws.Range("N11").Value = -3490
$ws.Range("H56").Value = 12333.333
$ws.Range("I56").Value = 12333.333
$ws.Range("K56").Value = 12333.333
$ws.Range("M56").Value = -11803.333
$ws.Range("H63").Value = 5252.75
$ws.Range("I63").Value = 505.5
$ws.Range("K63").Value = 1516.5
$ws.Range("M63").Value = -767.5
$ws.Range("H66").Value = 5252.75
$ws.Range("I66").Value = 505.5
$ws.Range("K66").Value = 4549.5
$ws.Range("M66").Value = -805.5
$ws.Range("H75").Value = 4032.8333
$ws.Range("J75").Value = 3039.4
$ws.Range("L75").Value = 9118.200000000001
$ws.Range("N75").Value = -11114.2
$ws.Range("H78").Value = 4032.8333
$ws.Range("J78").Value = 3039.4
$ws.Range("L78").Value = 27354.6
$ws.Range("N78").Value = -37338.60000000001
$ws.Range("H80").Value = 11999.25
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 11999.25
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 35997.75
$ws.Range("N80").Value = -37869.75
$ws.Range("M80").ClearContents()
$ws.Range("H81").Value = 129421.625
$ws.Range("I81").Value = 2964.6667
$ws.Range("J81").Value = 205295.8
$ws.Range("K81").Value = 8894.000100000001
$ws.Range("L81").Value = 615887.3999999999
$ws.Range("M81").Value = -7771.000100000001
$ws.Range("N81").Value = -618133.3999999999
$ws.Range("H83").Value = 11999.25
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 11999.25
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 107993.25
$ws.Range("N83").Value = -117353.25
$ws.Range("M83").ClearContents()
$ws.Range("H84").Value = 129421.625
$ws.Range("I84").Value = 2964.6667
$ws.Range("J84").Value = 205295.8
$ws.Range("K84").Value = 26682.0003
$ws.Range("L84").Value = 1847662.2
$ws.Range("M84").Value = -21066.0003
$ws.Range("N84").Value = -1858894.2
$ws.Range("I87").Value = 5205
$ws.Range("J87").Value = 9006
$ws.Range("K87").Value = 15615
$ws.Range("L87").Value = 27018
$ws.Range("M87").Value = -14367
$ws.Range("N87").Value = -29514
$ws.Range("I90").Value = 5205
$ws.Range("J90").Value = 9006
$ws.Range("K90").Value = 46845
$ws.Range("L90").Value = 81054
$ws.Range("M90").Value = -40605
$ws.Range("N90").Value = -93534
$ws.Range("H129").Value = 1768
$ws.Range("I129").Value = 1288
$ws.Range("J129").Value = 1848
$ws.Range("K129").Value = 3864
$ws.Range("L129").Value = 5544
$ws.Range("M129").Value = 1136
$ws.Range("N129").Value = -15544
$ws.Range("H131").Value = 2933.8276
$ws.Range("I131").Value = 1009.6667
$ws.Range("J131").Value = 3155.8462
$ws.Range("K131").Value = 3029.0001
$ws.Range("L131").Value = 9467.5386
$ws.Range("M131").Value = 2010.9999
$ws.Range("N131").Value = -19547.5386
$ws.Range("H137").Value = 1206
$ws.Range("I137").Value = 1206
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3618
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 1482
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2950
$ws.Range("I7").Value = 2950
$ws.Range("K7").Value = 2950
$ws.Range("M7").Value = -2838
$ws.Range("H82").Value = 2899.875
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2899.875
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H122").Value = 7199.8
$ws.Range("I122").Value = 7199.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 21599.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -19149.4
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2950
$ws.Range("I126").Value = 2950
$ws.Range("K126").Value = 8850
$ws.Range("M126").Value = -6380
